$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5, shifting existing rows 5-11 down to 6-12
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with values, following the same pattern
# as the other data rows (A=scaler name, C left as empty text, D/E fixed).
$ws.Cells.Item(5, 1).Value = "StandardScaler"
$ws.Cells.Item(5, 2).Value = 100
$ws.Cells.Item(5, 3).Value = "'"
$ws.Cells.Item(5, 3).Style = "Normal"
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 10000
$ws.Cells.Item(5, 7).Value = 2.302295207977295
$ws.Cells.Item(5, 8).Value = 0.00003260567561217538
